$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1881918819188192
$ws.Range("C2").Value = 0.5830258302583026
$ws.Range("J2").Value = 0.01476014760147601
$ws.Range("P2").Value = 0.1291512915129151
$ws.Range("S2").Value = 0.08487084870848709
$ws.Range("B3").Value = 0.005988023952095809
$ws.Range("C3").Value = 0.05389221556886228
$ws.Range("J3").Value = 0.03592814371257485
$ws.Range("P3").Value = 0.7245508982035929
$ws.Range("S3").Value = 0.1796407185628743
$ws.Range("J4").Value = 0.07843137254901961
$ws.Range("O4").Value = 0.0196078431372549
$ws.Range("P4").Value = 0.5882352941176471
$ws.Range("S4").Value = 0.3137254901960784
$ws.Range("B6").Value = 0.06542056074766354
$ws.Range("D6").Value = 0.009345794392523364
$ws.Range("F6").Value = 0.05607476635514019
$ws.Range("J6").Value = 0.2336448598130841
$ws.Range("O6").Value = 0.02803738317757009
$ws.Range("Q6").Value = 0.1635514018691589
$ws.Range("R6").Value = 0.03271028037383177
$ws.Range("S6").Value = 0.411214953271028
$ws.Range("B7").Value = 0.075
$ws.Range("D7").Value = 0.015
$ws.Range("F7").Value = 0.045
$ws.Range("J7").Value = 0.125
$ws.Range("O7").Value = 0.005
$ws.Range("Q7").Value = 0.15
$ws.Range("R7").Value = 0.06
$ws.Range("S7").Value = 0.525
$ws.Range("B8").Value = 0.07154742096505824
$ws.Range("D8").Value = 0.01996672212978369
$ws.Range("E8").Value = 0.001663893510815308
$ws.Range("F8").Value = 0.05324459234608985
$ws.Range("J8").Value = 0.129783693843594
$ws.Range("O8").Value = 0.021630615640599
$ws.Range("Q8").Value = 0.1447587354409318
$ws.Range("R8").Value = 0.09151414309484193
$ws.Range("S8").Value = 0.4658901830282862
$ws.Range("B9").Value = 0.1196581196581197
$ws.Range("D9").Value = 0.0170940170940171
$ws.Range("F9").Value = 0.04273504273504274
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("O9").Value = 0.008547008547008548
$ws.Range("Q9").Value = 0.1709401709401709
$ws.Range("R9").Value = 0.08974358974358974
$ws.Range("S9").Value = 0.4401709401709402
$ws.Range("B10").Value = 0.0891089108910891
$ws.Range("D10").Value = 0.02437166793602437
$ws.Range("F10").Value = 0.06626047220106626
$ws.Range("J10").Value = 0.1081492764661082
$ws.Range("O10").Value = 0.01904036557501904
$ws.Range("Q10").Value = 0.2254379284082254
$ws.Range("R10").Value = 0.0753998476770754
$ws.Range("S10").Value = 0.3922315308453923
$ws.Range("G11").Value = 0.124203821656051
$ws.Range("J11").Value = 0.07961783439490445
$ws.Range("K11").Value = 0.2292993630573248
$ws.Range("L11").Value = 0.5509554140127388
$ws.Range("S11").Value = 0.01592356687898089
$ws.Range("G12").Value = 0.7569060773480663
$ws.Range("J12").Value = 0.1823204419889503
$ws.Range("K12").Value = 0.01104972375690608
$ws.Range("L12").Value = 0.02209944751381215
$ws.Range("S12").Value = 0.02762430939226519
$ws.Range("G13").Value = 0.7659574468085106
$ws.Range("J13").Value = 0.2127659574468085
$ws.Range("S13").Value = 0.02127659574468085
$ws.Range("F15").Value = 0.008658008658008658
$ws.Range("H15").Value = 0.2207792207792208
$ws.Range("I15").Value = 0.08225108225108226
$ws.Range("J15").Value = 0.2857142857142857
$ws.Range("K15").Value = 0.06493506493506493
$ws.Range("M15").Value = 0.008658008658008658
$ws.Range("O15").Value = 0.05627705627705628
$ws.Range("S15").Value = 0.2727272727272727
$ws.Range("F16").Value = 0.02162162162162162
$ws.Range("H16").Value = 0.227027027027027
$ws.Range("I16").Value = 0.1027027027027027
$ws.Range("J16").Value = 0.345945945945946
$ws.Range("K16").Value = 0.1027027027027027
$ws.Range("M16").Value = 0.02702702702702703
$ws.Range("O16").Value = 0.04324324324324325
$ws.Range("S16").Value = 0.1297297297297297
$ws.Range("F17").Value = 0.01659751037344398
$ws.Range("H17").Value = 0.2448132780082987
$ws.Range("I17").Value = 0.07883817427385892
$ws.Range("J17").Value = 0.3983402489626556
$ws.Range("K17").Value = 0.08506224066390042
$ws.Range("M17").Value = 0.01244813278008299
$ws.Range("O17").Value = 0.06224066390041494
$ws.Range("S17").Value = 0.1016597510373444
$ws.Range("F18").Value = 0.02564102564102564
$ws.Range("H18").Value = 0.2051282051282051
$ws.Range("I18").Value = 0.1128205128205128
$ws.Range("J18").Value = 0.4153846153846154
$ws.Range("K18").Value = 0.1025641025641026
$ws.Range("M18").Value = 0.02051282051282051
$ws.Range("O18").Value = 0.03589743589743589
$ws.Range("S18").Value = 0.08205128205128205
$ws.Range("F19").Value = 0.01071667782987274
$ws.Range("H19").Value = 0.2344273275284662
$ws.Range("I19").Value = 0.09109176155391828
$ws.Range("J19").Value = 0.3556597454789016
$ws.Range("K19").Value = 0.0971198928332217
$ws.Range("M19").Value = 0.02009377093101139
$ws.Range("N19").Value = 0.0006697923643670462
$ws.Range("O19").Value = 0.06296048225050234
$ws.Range("S19").Value = 0.1272605492297388
